$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price/volume columns to Text format so numeric-looking
# strings (e.g. "22.70", "1.001") are preserved verbatim instead of
# being normalized into numbers by Excel's automatic type coercion.
$ws.Range("D2:E51").NumberFormat = "@"

$updates = @(
    @{ Cell = "D2"; Value = "29.311.71" }
    @{ Cell = "E2"; Value = "  +0.61%  " }
    @{ Cell = "D3"; Value = "1.934.42" }
    @{ Cell = "E3"; Value = "  +1.27%  " }
    @{ Cell = "D4"; Value = "1.003" }
    @{ Cell = "E4"; Value = "  +0.27%  " }
    @{ Cell = "D5"; Value = "325.57" }
    @{ Cell = "E5"; Value = "  +0.24%  " }
    @{ Cell = "D6"; Value = "1.001" }
    @{ Cell = "E6"; Value = "  +0.33%  " }
    @{ Cell = "D7"; Value = "0.4628" }
    @{ Cell = "E7"; Value = "  +0.43%  " }
    @{ Cell = "D8"; Value = "0.3877" }
    @{ Cell = "E8"; Value = "  -0.27%  " }
    @{ Cell = "D9"; Value = "45.96" }
    @{ Cell = "E9"; Value = "  +0.81%  " }
    @{ Cell = "D10"; Value = "0.07828" }
    @{ Cell = "E10"; Value = "  -0.38%  " }
    @{ Cell = "D11"; Value = "0.9748" }
    @{ Cell = "E11"; Value = "  -1.72%  " }
    @{ Cell = "D12"; Value = "22.70" }
    @{ Cell = "E12"; Value = "  +3.11%  " }
    @{ Cell = "D13"; Value = "1.932.68" }
    @{ Cell = "E13"; Value = "  +2.30%  " }
    @{ Cell = "D14"; Value = "7.087" }
    @{ Cell = "E14"; Value = "  +0.61%  " }
    @{ Cell = "E15"; Value = "  +0.52%  " }
    @{ Cell = "D16"; Value = "0.07081" }
    @{ Cell = "E16"; Value = "  +0.67%  " }
    @{ Cell = "D17"; Value = "86.78" }
    @{ Cell = "E17"; Value = "  -1.28%  " }
    @{ Cell = "D18"; Value = "1.004" }
    @{ Cell = "E18"; Value = "  +0.33%  " }
    @{ Cell = "D19"; Value = "0.000009739" }
    @{ Cell = "E19"; Value = "  -1.92%  " }
    @{ Cell = "D20"; Value = "17.02" }
    @{ Cell = "E20"; Value = "  -0.22%  " }
    @{ Cell = "D21"; Value = "1.001" }
    @{ Cell = "E21"; Value = "  +0.28%  " }
    @{ Cell = "D22"; Value = "29.289.03" }
    @{ Cell = "E22"; Value = "  +0.51%  " }
    @{ Cell = "D23"; Value = "5.481" }
    @{ Cell = "E23"; Value = "  +2.99%  " }
    @{ Cell = "D24"; Value = "11.08" }
    @{ Cell = "E24"; Value = "  -0.38%  " }
    @{ Cell = "D25"; Value = "2.168.25" }
    @{ Cell = "E25"; Value = "  +1.94%  " }
    @{ Cell = "D26"; Value = "2.094" }
    @{ Cell = "E26"; Value = "  +0.43%  " }
    @{ Cell = "D27"; Value = "157.93" }
    @{ Cell = "E27"; Value = "  +1.25%  " }
    @{ Cell = "D28"; Value = "19.41" }
    @{ Cell = "E28"; Value = "  -0.22%  " }
    @{ Cell = "D29"; Value = "5.772" }
    @{ Cell = "E29"; Value = "  -2.36%  " }
    @{ Cell = "D30"; Value = "118.89" }
    @{ Cell = "E30"; Value = "  +0.20%  " }
    @{ Cell = "D31"; Value = "1.838" }
    @{ Cell = "E31"; Value = "  -1.64%  " }
    @{ Cell = "D32"; Value = "0.09336" }
    @{ Cell = "E32"; Value = "  +0.01%  " }
    @{ Cell = "D33"; Value = "0.8659" }
    @{ Cell = "E33"; Value = "  -3.14%  " }
    @{ Cell = "D34"; Value = "5.167" }
    @{ Cell = "E34"; Value = "  -0.99%  " }
    @{ Cell = "D35"; Value = "1.301" }
    @{ Cell = "E35"; Value = "  -1.49%  " }
    @{ Cell = "E36"; Value = "  -1.84%  " }
    @{ Cell = "D37"; Value = "0.05786" }
    @{ Cell = "E37"; Value = "  +0.13%  " }
    @{ Cell = "D38"; Value = "1.158" }
    @{ Cell = "E38"; Value = "  -0.81%  " }
    @{ Cell = "D39"; Value = "0.02080" }
    @{ Cell = "E39"; Value = "  -0.27%  " }
    @{ Cell = "D40"; Value = "7.644" }
    @{ Cell = "E40"; Value = "  +0.09%  " }
    @{ Cell = "D41"; Value = "0.5659" }
    @{ Cell = "E41"; Value = "  -0.67%  " }
    @{ Cell = "B42"; Value = "Algorand" }
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo" }
    @{ Cell = "D42"; Value = "0.1779" }
    @{ Cell = "E42"; Value = "  -1.75%  " }
    @{ Cell = "B43"; Value = "PEPE" }
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe" }
    @{ Cell = "D43"; Value = "0.000003050" }
    @{ Cell = "E43"; Value = "  +2.01%  " }
    @{ Cell = "D44"; Value = "9.407" }
    @{ Cell = "E44"; Value = "  -3.35%  " }
    @{ Cell = "D45"; Value = "2.711" }
    @{ Cell = "E45"; Value = "  +6.77%  " }
    @{ Cell = "D46"; Value = "0.5272" }
    @{ Cell = "E46"; Value = "  -1.51%  " }
    @{ Cell = "D47"; Value = "11.50" }
    @{ Cell = "E47"; Value = "  -3.09%  " }
    @{ Cell = "D48"; Value = "0.06861" }
    @{ Cell = "E48"; Value = "  -1.56%  " }
    @{ Cell = "D49"; Value = "2.076" }
    @{ Cell = "E49"; Value = "  -4.71%  " }
    @{ Cell = "E50"; Value = "  -1.72%  " }
    @{ Cell = "D51"; Value = "111.40" }
    @{ Cell = "E51"; Value = "  -1.43%  " }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
